# Update automàtic: dades i banners [2026-02-14 23:20]
#
# Refreshes the per-station extraction timestamps (col E) and the observed
# weather readings that moved between the 22:49 and 23:20 extraction runs.
#
# Note on column H (HUMITAT_MITJANA_DIA, percentages stored as plain text
# like "71%"): assigning that text straight to .Value makes Excel's normal
# type-inference reinterpret it as a numeric percentage (0.71, formatted
# "#,##0%"), which would silently change the cell's style/number format.
# To keep it as literal text with the original "General" style untouched,
# the value is entered with a leading apostrophe (forces text entry, same
# as typing it in Excel) and then the original cell formatting is restored
# by pasting (formats only) from C2 - a same-style (s="3"), never-edited
# reference cell - onto the cell that was just overwritten.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-14 23:18:26"
$ws.Range("K2").Value = "1.8 MJ/m2"
$ws.Range("N2").Value = "-4.5 °C 22:45 TU"
$ws.Range("O2").Value = "-1.4 °C"
$ws.Range("E3").Value = "2026-02-14 23:18:28"
$ws.Range("N3").Value = "-9.1 °C 22:49 TU"
$ws.Range("O3").Value = "-5.5 °C"
$ws.Range("E4").Value = "2026-02-14 23:18:30"
$ws.Range("J4").Value = "998.6 hPa"
$ws.Range("N4").Value = "3.7 °C 22:50 TU"
$ws.Range("O4").Value = "10.4 °C"
$ws.Range("E5").Value = "2026-02-14 23:18:33"
$ws.Range("N5").Value = "-8.5 °C 22:59 TU"
$ws.Range("E6").Value = "2026-02-14 23:18:35"
$ws.Range("H6").Value = "'71%"
$ws.Range("C2").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("J6").Value = "998.6 hPa"
$ws.Range("E7").Value = "2026-02-14 23:18:38"
$ws.Range("H7").Value = "'49%"
$ws.Range("C2").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("J7").Value = "998.9 hPa"
$ws.Range("E8").Value = "2026-02-14 23:18:40"
$ws.Range("H8").Value = "'59%"
$ws.Range("C2").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("J8").Value = "998.7 hPa"
$ws.Range("E9").Value = "2026-02-14 23:18:43"
$ws.Range("L9").Value = "74.2 km/h - 0º 22:41 TU"
$ws.Range("E10").Value = "2026-02-14 23:18:46"
$ws.Range("H10").Value = "'73%"
$ws.Range("C2").Copy()
$ws.Range("H10").PasteSpecial(-4122)
$ws.Range("E11").Value = "2026-02-14 23:18:48"
$ws.Range("E12").Value = "2026-02-14 23:18:51"
$ws.Range("O12").Value = "11.9 °C"
$ws.Range("E13").Value = "2026-02-14 23:18:53"
$ws.Range("H13").Value = "'65%"
$ws.Range("C2").Copy()
$ws.Range("H13").PasteSpecial(-4122)
$ws.Range("J13").Value = "1001.6 hPa"
$ws.Range("L13").Value = "84.2 km/h - 16º 22:55 TU"
$ws.Range("E14").Value = "2026-02-14 23:18:55"
$ws.Range("H14").Value = "'49%"
$ws.Range("C2").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("O14").Value = "13.2 °C"
$ws.Range("E15").Value = "2026-02-14 23:18:58"
$ws.Range("E16").Value = "2026-02-14 23:19:00"
$ws.Range("H16").Value = "'73%"
$ws.Range("C2").Copy()
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("M16").Value = "-3.8 °C 22:59 TU"
$ws.Range("E17").Value = "2026-02-14 23:19:03"
$ws.Range("H17").Value = "'65%"
$ws.Range("C2").Copy()
$ws.Range("H17").PasteSpecial(-4122)
$ws.Range("E18").Value = "2026-02-14 23:19:05"
$ws.Range("H18").Value = "'71%"
$ws.Range("C2").Copy()
$ws.Range("H18").PasteSpecial(-4122)
$ws.Range("J18").Value = "998.8 hPa"
$ws.Range("L18").Value = "31.7 km/h - 33º 22:46 TU"
$ws.Range("E19").Value = "2026-02-14 23:19:08"
$ws.Range("O19").Value = "5.6 °C"
$ws.Range("E20").Value = "2026-02-14 23:19:10"
$ws.Range("N20").Value = "-8.9 °C 22:46 TU"
$ws.Range("O20").Value = "-5.7 °C"
$ws.Range("E21").Value = "2026-02-14 23:19:12"
$ws.Range("H21").Value = "'66%"
$ws.Range("C2").Copy()
$ws.Range("H21").PasteSpecial(-4122)
$ws.Range("J21").Value = "1001.3 hPa"
$ws.Range("E22").Value = "2026-02-14 23:19:15"
$ws.Range("H22").Value = "'83%"
$ws.Range("C2").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("E23").Value = "2026-02-14 23:19:17"
$ws.Range("N23").Value = "-9.4 °C 22:59 TU"
$ws.Range("E24").Value = "2026-02-14 23:19:20"
$ws.Range("J24").Value = "1003.0 hPa"
$ws.Range("O24").Value = "9.2 °C"
$ws.Range("E25").Value = "2026-02-14 23:19:22"
$ws.Range("H25").Value = "'83%"
$ws.Range("C2").Copy()
$ws.Range("H25").PasteSpecial(-4122)
$ws.Range("I25").Value = "20.9 mm"
$ws.Range("O25").Value = "-5.0 °C"
$ws.Range("E26").Value = "2026-02-14 23:19:25"
$ws.Range("E27").Value = "2026-02-14 23:19:27"
$ws.Range("H27").Value = "'74%"
$ws.Range("C2").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("L27").Value = "61.6 km/h - 305º 22:43 TU"
$ws.Range("E28").Value = "2026-02-14 23:19:30"
$ws.Range("J28").Value = "998.5 hPa"
$ws.Range("O28").Value = "9.4 °C"
$ws.Range("E29").Value = "2026-02-14 23:19:32"
$ws.Range("H29").Value = "'61%"
$ws.Range("C2").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("O29").Value = "11.4 °C"
$ws.Range("E30").Value = "2026-02-14 23:19:35"
$ws.Range("J30").Value = "998.5 hPa"
$ws.Range("E31").Value = "2026-02-14 23:19:37"
$ws.Range("J31").Value = "997.8 hPa"
$ws.Range("N31").Value = "6.9 °C 22:38 TU"
$ws.Range("E32").Value = "2026-02-14 23:19:40"
$ws.Range("H32").Value = "'85%"
$ws.Range("C2").Copy()
$ws.Range("H32").PasteSpecial(-4122)
$ws.Range("N32").Value = "1.5 °C 22:59 TU"
$ws.Range("O32").Value = "4.0 °C"
$ws.Range("E33").Value = "2026-02-14 23:19:42"
$ws.Range("H33").Value = "'61%"
$ws.Range("C2").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("J33").Value = "1000.9 hPa"
$ws.Range("E34").Value = "2026-02-14 23:19:45"
$ws.Range("N34").Value = "-5.5 °C 22:54 TU"
$ws.Range("O34").Value = "-2.6 °C"
$ws.Range("E35").Value = "2026-02-14 23:19:47"
$ws.Range("J35").Value = "1005.4 hPa"
$ws.Range("N35").Value = "1.1 °C 22:58 TU"
$ws.Range("E36").Value = "2026-02-14 23:19:50"
$ws.Range("J36").Value = "999.3 hPa"
$ws.Range("E37").Value = "2026-02-14 23:19:52"
$ws.Range("H37").Value = "'61%"
$ws.Range("C2").Copy()
$ws.Range("H37").PasteSpecial(-4122)
$ws.Range("J37").Value = "999.5 hPa"
$ws.Range("E38").Value = "2026-02-14 23:19:55"
$ws.Range("O38").Value = "9.9 °C"
$ws.Range("E39").Value = "2026-02-14 23:19:57"
$ws.Range("E40").Value = "2026-02-14 23:20:00"
$ws.Range("H40").Value = "'63%"
$ws.Range("C2").Copy()
$ws.Range("H40").PasteSpecial(-4122)
$ws.Range("J40").Value = "1002.0 hPa"
$ws.Range("E41").Value = "2026-02-14 23:20:02"
$ws.Range("J41").Value = "1000.7 hPa"
$ws.Range("N41").Value = "10.1 °C 22:58 TU"
$ws.Range("E42").Value = "2026-02-14 23:20:05"
$ws.Range("E43").Value = "2026-02-14 23:20:07"
$ws.Range("H43").Value = "'62%"
$ws.Range("C2").Copy()
$ws.Range("H43").PasteSpecial(-4122)
$ws.Range("K43").Value = "13.3 MJ/m2"
$ws.Range("O43").Value = "9.1 °C"
$ws.Range("E44").Value = "2026-02-14 23:20:09"
$ws.Range("H44").Value = "'91%"
$ws.Range("C2").Copy()
$ws.Range("H44").PasteSpecial(-4122)
$ws.Range("I44").Value = "37.9 mm"
$ws.Range("N44").Value = "-8.8 °C 22:57 TU"
$ws.Range("O44").Value = "-5.7 °C"
$ws.Range("E45").Value = "2026-02-14 23:20:12"
$ws.Range("J45").Value = "1008.2 hPa"
$ws.Range("N45").Value = "-0.2 °C 22:56 TU"
$ws.Range("O45").Value = "2.7 °C"
$ws.Range("E46").Value = "2026-02-14 23:20:14"

$excel.CutCopyMode = $false

